$d = $word.ActiveDocument

# --- 1. Rework the "defendant response deadline" paragraph into the
#        conditional cs_/else/es_ block (registered vs unregistered org). ---
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*If you do*not respond before*County Court Judgment*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:after="40"/></w:pPr><w:r><w:t>&lt;&lt;cs_{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>respondentsOrgRegistered</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=' + [char]0x2019 + 'No' + [char]0x2019 + '}&gt;&gt;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:after="40"/></w:pPr><w:r><w:t>If you do not respond, a county court judgment could be issued.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:after="40"/></w:pPr><w:r><w:t>&lt;&lt;else&gt;&gt;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:after="40"/></w:pPr><w:r><w:t xml:space="preserve">If you do not respond before </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>&lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>defendantResponseDeadlineDate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>&gt;&gt;</w:t></w:r><w:r><w:t>, you could get a County Court Judgment (CCJ) made against you.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:after="40"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>&lt;&lt;es_&gt;&gt;</w:t></w:r></w:p>'

    $target.Range.InsertXML($newXml)
}

# --- 2. styles.xml: the "Default Paragraph Font" character style is no
#        longer semi-hidden. ---
$styles = $d.Styles
$dpf = $styles.Item("Default Paragraph Font")
$dpf.Font.Hidden = $false
